$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.672.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.312.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.313.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.859.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.315.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.765.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("E30").Value = "  -4.41%  "
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0744"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0399"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.123.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "432.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E42").Value = "  +8.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E46").Value = "  +3.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("E51").Value = "  -0.52%  "
